$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J6").Value = 3
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.73
$ws.Range("AE6").Value = 19
$ws.Range("G7").Value = 2.45
$ws.Range("H7").Value = 2.88
$ws.Range("I7").Value = 3.4
$ws.Range("J7").Value = 3.4
$ws.Range("K7").Value = 1.8
$ws.Range("L7").Value = 4.33
$ws.Range("M7").Value = 1.14
$ws.Range("N7").Value = 5.5
$ws.Range("O7").Value = 1.67
$ws.Range("P7").Value = 2.1
$ws.Range("Q7").Value = 3.1
$ws.Range("R7").Value = 1.36
$ws.Range("S7").Value = 5.2
$ws.Range("T7").Value = 1.16
$ws.Range("U7").Value = 6.5
$ws.Range("V7").Value = 1.11
$ws.Range("Y7").Value = 2.38
$ws.Range("Z7").Value = 1.53
$ws.Range("AB7").Value = 10
$ws.Range("AD7").Value = 23
$ws.Range("AE7").Value = 26
$ws.Range("AJ7").Value = 101
$ws.Range("AL7").Value = 7
$ws.Range("AN7").Value = 13
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 4.33
$ws.Range("J8").Value = 3.1
$ws.Range("L8").Value = 5.5
$ws.Range("M8").Value = 1.2
$ws.Range("N8").Value = 4.33
$ws.Range("O8").Value = 1.83
$ws.Range("P8").Value = 1.83
$ws.Range("Q8").Value = 4
$ws.Range("R8").Value = 1.25
$ws.Range("U8").Value = 10
$ws.Range("V8").Value = 1.06
$ws.Range("W8").Value = 1.85
$ws.Range("X8").Value = 1.95
$ws.Range("Y8").Value = 3
$ws.Range("Z8").Value = 1.36
$ws.Range("AA8").Value = 4.33
$ws.Range("AB8").Value = 8
$ws.Range("AH8").Value = 6.5
$ws.Range("AN8").Value = 19
$ws.Range("AQ8").Value = 81
$ws.Range("G10").Value = 1.8
$ws.Range("I10").Value = 5.5
$ws.Range("J10").Value = 2.6
$ws.Range("K10").Value = 1.91
$ws.Range("L10").Value = 6
$ws.Range("M10").Value = 1.13
$ws.Range("N10").Value = 6
$ws.Range("O10").Value = 1.53
$ws.Range("P10").Value = 2.38
$ws.Range("Q10").Value = 2.7
$ws.Range("R10").Value = 1.44
$ws.Range("S10").Value = 4.4
$ws.Range("T10").Value = 1.2
$ws.Range("U10").Value = 5.5
$ws.Range("V10").Value = 1.14
$ws.Range("W10").Value = 1.62
$ws.Range("X10").Value = 2.2
$ws.Range("Y10").Value = 2.38
$ws.Range("Z10").Value = 1.53
$ws.Range("AA10").Value = 5
$ws.Range("AB10").Value = 7
$ws.Range("AC10").Value = 10
$ws.Range("AD10").Value = 13
$ws.Range("AG10").Value = 5.5
$ws.Range("AI10").Value = 21
$ws.Range("AM10").Value = 26
$ws.Range("AN10").Value = 19
$ws.Range("AP10").Value = 51
$ws.Range("AQ10").Value = 67
$ws.Range("AR10").Value = 2.1
$ws.Range("AS10").Value = 1.78
$ws.Range("Q11").Value = 2.6
$ws.Range("R11").Value = 1.48
$ws.Range("S11").Value = 4.1
$ws.Range("T11").Value = 1.22
$ws.Range("U11").Value = 5.5
$ws.Range("V11").Value = 1.14
$ws.Range("AR11").Value = 1.98
$ws.Range("AS11").Value = 1.88
